# issue #5: stock data from json to db
# The "股票" (stock) sheet gains three new columns that mirror the
# fields added to the JSON->DB export: a "category" column (value
# "normal") inserted right after "property_category", and two columns
# appended at the end: "source_file" (value "tmp61ee1") and "index"
# (the original per-row id that lives in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 10

# Insert a new column before the current "date" column (col 9 / I) and
# turn it into "category" = "normal" for every data row.
$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# Append "source_file" and "index" columns (M, N) after legislator_id (L).
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
for ($r = 2; $r -le $lastRow; $r++) {
    $origIndex = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 13).Value = "tmp61ee1"
    $ws.Cells.Item($r, 14).Value = $origIndex
}

# Match the header formatting (bold + border, same style as the other
# header cells) for the two newly appended header cells.
$ws.Range("H1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
